# "status info, repeat dialog / update"
#
# The Skill Data XML-mapped table gains a new `skillMaxStack` attribute/
# column (9th column, "I"), the workbook's active tab moves from
# "NPC Dialog Data" to "Skill Data", and the new cell I15 becomes the
# live selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill Data")

# --- new "skillMaxStack" column (I) ------------------------------------
$ws.Range("I1").Value = "skillMaxStack"

# Row 2..15 skillMaxStack values, in sheet order. These come through the
# XML map as text (not numbers), so force each cell to a text format
# before assigning, otherwise Excel would coerce the numeric-looking
# strings into real numbers.
$skillMaxStack = @("1","1","1","1","1","1","1","1","1","3","5","5","3","3")
for ($i = 0; $i -lt $skillMaxStack.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("I" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $skillMaxStack[$i]
}

# Best-effort cosmetic fit for the new column's width.
$ws.Columns.Item(9).AutoFit() | Out-Null

# --- active tab / selection ---------------------------------------------
# Previously "NPC Dialog Data" (activeTab=2) was active with a selection
# on its own sheet; now "Skill Data" (activeTab=1) becomes active, with
# the selection resting on the newly-entered I15 cell.
$ws.Activate() | Out-Null
$ws.Range("I15").Select() | Out-Null
